$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.037278529888334
$ws.Cells.Item(2, 4).Value = 1.043896643453357
$ws.Cells.Item(2, 5).Value = 1.036084749279061
$ws.Cells.Item(2, 6).Value = 1.052647421909903
$ws.Cells.Item(2, 9).Value = 1.039039805171452
$ws.Cells.Item(2, 10).Value = 1.042382269688147
$ws.Cells.Item(2, 11).Value = 1.046669217456616
$ws.Cells.Item(2, 12).Value = 1.038879501990603
$ws.Cells.Item(2, 13).Value = 1.055395575432976
$ws.Cells.Item(2, 14).Value = 1.018030650475702

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038236254006643
$ws.Cells.Item(3, 4).Value = 1.0446549255685
$ws.Cells.Item(3, 5).Value = 1.03689961089894
$ws.Cells.Item(3, 6).Value = 1.05363866180943
$ws.Cells.Item(3, 9).Value = 1.039280819455312
$ws.Cells.Item(3, 10).Value = 1.042984265356505
$ws.Cells.Item(3, 11).Value = 1.04723876573907
$ws.Cells.Item(3, 12).Value = 1.039503877975598
$ws.Cells.Item(3, 13).Value = 1.05619924513287
$ws.Cells.Item(3, 14).Value = 1.018233852474666

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.038856374692977
$ws.Cells.Item(4, 4).Value = 1.045145957218835
$ws.Cells.Item(4, 5).Value = 1.037427589661962
$ws.Cells.Item(4, 6).Value = 1.054280957339708
$ws.Cells.Item(4, 9).Value = 1.039435790063998
$ws.Cells.Item(4, 10).Value = 1.043373591488271
$ws.Cells.Item(4, 11).Value = 1.047607010653035
$ws.Cells.Item(4, 12).Value = 1.03990794948473
$ws.Cells.Item(4, 13).Value = 1.056719553776026
$ws.Cells.Item(4, 14).Value = 1.018365182562633

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.039117170542377
$ws.Cells.Item(5, 4).Value = 1.045352474857227
$ws.Cells.Item(5, 5).Value = 1.037649720378461
$ws.Cells.Item(5, 6).Value = 1.054551191948855
$ws.Cells.Item(5, 9).Value = 1.039500704189225
$ws.Cells.Item(5, 10).Value = 1.0435372141769
$ws.Cells.Item(5, 11).Value = 1.047761750282291
$ws.Cells.Item(5, 12).Value = 1.040077834204238
$ws.Cells.Item(5, 13).Value = 1.056938357653175
$ws.Cells.Item(5, 14).Value = 1.018420356208726

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.039160964983603
$ws.Cells.Item(6, 4).Value = 1.045387155198258
$ws.Cells.Item(6, 5).Value = 1.037687026944672
$ws.Cells.Item(6, 6).Value = 1.054596578018245
$ws.Cells.Item(6, 9).Value = 1.039511589731483
$ws.Cells.Item(6, 10).Value = 1.04356468418822
$ws.Cells.Item(6, 11).Value = 1.047787727585086
$ws.Cells.Item(6, 12).Value = 1.040106359334347
$ws.Cells.Item(6, 13).Value = 1.056975099633863
$ws.Cells.Item(6, 14).Value = 1.018429617895589

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.038859859079327
$ws.Cells.Item(7, 4).Value = 1.045148716373645
$ws.Cells.Item(7, 5).Value = 1.037430557123156
$ws.Cells.Item(7, 6).Value = 1.054284567390079
$ws.Cells.Item(7, 9).Value = 1.039436658375745
$ws.Cells.Item(7, 10).Value = 1.043375778020254
$ws.Cells.Item(7, 11).Value = 1.047609078569788
$ws.Cells.Item(7, 12).Value = 1.039910219442545
$ws.Cells.Item(7, 13).Value = 1.056722477185636
$ws.Cells.Item(7, 14).Value = 1.018365919943493

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037602112018745
$ws.Cells.Item(8, 4).Value = 1.044152830596516
$ws.Cells.Item(8, 5).Value = 1.036359988066714
$ws.Cells.Item(8, 6).Value = 1.052982229625955
$ws.Cells.Item(8, 9).Value = 1.039121460061962
$ws.Cells.Item(8, 10).Value = 1.042585758936189
$ws.Cells.Item(8, 11).Value = 1.046861758665472
$ws.Cells.Item(8, 12).Value = 1.039090500096389
$ws.Cells.Item(8, 13).Value = 1.055667120634471
$ws.Cells.Item(8, 14).Value = 1.018099355454323

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.03538897144332
$ws.Cells.Item(9, 4).Value = 1.042400858941688
$ws.Cells.Item(9, 5).Value = 1.034478992182175
$ws.Cells.Item(9, 6).Value = 1.050694266559379
$ws.Cells.Item(9, 9).Value = 1.038558542829174
$ws.Cells.Item(9, 10).Value = 1.041192109963547
$ws.Cells.Item(9, 11).Value = 1.045542696708216
$ws.Cells.Item(9, 12).Value = 1.037646539679249
$ws.Cells.Item(9, 13).Value = 1.053809643238017
$ws.Cells.Item(9, 14).Value = 1.01762846246099

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033915724323855
$ws.Cells.Item(10, 4).Value = 1.041234904030708
$ws.Cells.Item(10, 5).Value = 1.033228755124746
$ws.Cells.Item(10, 6).Value = 1.049173681501981
$ws.Cells.Item(10, 9).Value = 1.038178252325725
$ws.Cells.Item(10, 10).Value = 1.040262029554822
$ws.Cells.Item(10, 11).Value = 1.044661902587894
$ws.Cells.Item(10, 12).Value = 1.036684284338336
$ws.Cells.Item(10, 13).Value = 1.0525728630632
$ws.Cells.Item(10, 14).Value = 1.017313767307285

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033278318807507
$ws.Cells.Item(11, 4).Value = 1.040730528243643
$ws.Cells.Item(11, 5).Value = 1.032688295775276
$ws.Cells.Item(11, 6).Value = 1.048516386047764
$ws.Cells.Item(11, 9).Value = 1.038012399351212
$ws.Cells.Item(11, 10).Value = 1.039859072158302
$ws.Cells.Item(11, 11).Value = 1.044280183079952
$ws.Cells.Item(11, 12).Value = 1.036267719944407
$ws.Cells.Item(11, 13).Value = 1.052037701133741
$ws.Cells.Item(11, 14).Value = 1.01717732351935

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.033041636893176
$ws.Cells.Item(12, 4).Value = 1.040543255471946
$ws.Cells.Item(12, 5).Value = 1.032487681758369
$ws.Cells.Item(12, 6).Value = 1.048272407815987
$ws.Cells.Item(12, 9).Value = 1.037950616628536
$ws.Cells.Item(12, 10).Value = 1.03970936261741
$ws.Cells.Item(12, 11).Value = 1.044138346733439
$ws.Cells.Item(12, 12).Value = 1.036113005076226
$ws.Cells.Item(12, 13).Value = 1.051838975045857
$ws.Cells.Item(12, 14).Value = 1.017126615720902

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03309240238637
$ws.Cells.Item(13, 4).Value = 1.040583422711009
$ws.Cells.Item(13, 5).Value = 1.032530707943013
$ws.Cells.Item(13, 6).Value = 1.048324734238574
$ws.Cells.Item(13, 9).Value = 1.037963877259379
$ws.Cells.Item(13, 10).Value = 1.039741477312666
$ws.Cells.Item(13, 11).Value = 1.044168773302816
$ws.Cells.Item(13, 12).Value = 1.036146191216032
$ws.Cells.Item(13, 13).Value = 1.05188159987648
$ws.Cells.Item(13, 14).Value = 1.017137493909402

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033258752993427
$ws.Cells.Item(14, 4).Value = 1.04071504668412
$ws.Cells.Item(14, 5).Value = 1.032671710163493
$ws.Cells.Item(14, 6).Value = 1.048496215245086
$ws.Cells.Item(14, 9).Value = 1.038007295991352
$ws.Cells.Item(14, 10).Value = 1.039846697794045
$ws.Cells.Item(14, 11).Value = 1.044268459832475
$ws.Cells.Item(14, 12).Value = 1.036254930839785
$ws.Cells.Item(14, 13).Value = 1.052021273207264
$ws.Cells.Item(14, 14).Value = 1.017173132535496

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033361257566014
$ws.Cells.Item(15, 4).Value = 1.040796154512095
$ws.Cells.Item(15, 5).Value = 1.032758604428542
$ws.Cells.Item(15, 6).Value = 1.048601892987711
$ws.Cells.Item(15, 9).Value = 1.038034024192143
$ws.Cells.Item(15, 10).Value = 1.039911523226998
$ws.Cells.Item(15, 11).Value = 1.044329873561764
$ws.Cells.Item(15, 12).Value = 1.036321931017378
$ws.Cells.Item(15, 13).Value = 1.052107338122935
$ws.Cells.Item(15, 14).Value = 1.017195087168121

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033958037850206
$ws.Cells.Item(16, 4).Value = 1.041268388222959
$ws.Cells.Item(16, 5).Value = 1.03326464272897
$ws.Cells.Item(16, 6).Value = 1.049217327894983
$ws.Cells.Item(16, 9).Value = 1.038189234519217
$ws.Cells.Item(16, 10).Value = 1.040288767807329
$ws.Cells.Item(16, 11).Value = 1.04468722917587
$ws.Cells.Item(16, 12).Value = 1.036711932497186
$ws.Cells.Item(16, 13).Value = 1.052608387920555
$ws.Cells.Item(16, 14).Value = 1.017322818896674

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034332522136519
$ws.Cells.Item(17, 4).Value = 1.041564740074812
$ws.Cells.Item(17, 5).Value = 1.033582309403672
$ws.Cells.Item(17, 6).Value = 1.049603676708997
$ws.Cells.Item(17, 9).Value = 1.038286276995608
$ws.Cells.Item(17, 10).Value = 1.040525343303769
$ws.Cells.Item(17, 11).Value = 1.044911301188603
$ws.Cells.Item(17, 12).Value = 1.036956597128663
$ws.Cells.Item(17, 13).Value = 1.052922783523098
$ws.Cells.Item(17, 14).Value = 1.017402894035805

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034551002628969
$ws.Cells.Item(18, 4).Value = 1.041737644331968
$ws.Cells.Item(18, 5).Value = 1.033767685906066
$ws.Cells.Item(18, 6).Value = 1.049829136188021
$ws.Cells.Item(18, 9).Value = 1.03834276589639
$ws.Cells.Item(18, 10).Value = 1.04066331188629
$ws.Cells.Item(18, 11).Value = 1.045041966721957
$ws.Cells.Item(18, 12).Value = 1.0370993153144
$ws.Cells.Item(18, 13).Value = 1.053106200989911
$ws.Cells.Item(18, 14).Value = 1.017449583243414

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034625507314337
$ws.Cells.Item(19, 4).Value = 1.041796608191587
$ws.Cells.Item(19, 5).Value = 1.033830909206942
$ws.Cells.Item(19, 6).Value = 1.049906030456524
$ws.Cells.Item(19, 9).Value = 1.038362007745367
$ws.Cells.Item(19, 10).Value = 1.040710351868288
$ws.Cells.Item(19, 11).Value = 1.0450865148594
$ws.Cells.Item(19, 12).Value = 1.037147980102527
$ws.Cells.Item(19, 13).Value = 1.053168747653805
$ws.Cells.Item(19, 14).Value = 1.017465500117384

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.034292338324809
$ws.Cells.Item(20, 4).Value = 1.041532939431398
$ws.Cells.Item(20, 5).Value = 1.033548217778093
$ws.Cells.Item(20, 6).Value = 1.049562213880809
$ws.Cells.Item(20, 9).Value = 1.038275877079186
$ws.Cells.Item(20, 10).Value = 1.040499963251568
$ws.Cells.Item(20, 11).Value = 1.044887263673787
$ws.Cells.Item(20, 12).Value = 1.036930345945708
$ws.Cells.Item(20, 13).Value = 1.052889048137046
$ws.Cells.Item(20, 14).Value = 1.017394304514607

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033209764702688
$ws.Cells.Item(21, 4).Value = 1.040676284628905
$ws.Cells.Item(21, 5).Value = 1.032630184742098
$ws.Cells.Item(21, 6).Value = 1.048445713645329
$ws.Cells.Item(21, 9).Value = 1.03799451515448
$ws.Cells.Item(21, 10).Value = 1.039815713895503
$ws.Cells.Item(21, 11).Value = 1.044239105972492
$ws.Cells.Item(21, 12).Value = 1.036222909290402
$ws.Cells.Item(21, 13).Value = 1.051980141314501
$ws.Cells.Item(21, 14).Value = 1.017162638577122

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032529564405652
$ws.Cells.Item(22, 4).Value = 1.040138105414007
$ws.Cells.Item(22, 5).Value = 1.03205377214291
$ws.Cells.Item(22, 6).Value = 1.047744713028087
$ws.Cells.Item(22, 9).Value = 1.037816584754474
$ws.Cells.Item(22, 10).Value = 1.039385307351345
$ws.Cells.Item(22, 11).Value = 1.043831302003707
$ws.Cells.Item(22, 12).Value = 1.035778207038651
$ws.Cells.Item(22, 13).Value = 1.051409005105405
$ws.Cells.Item(22, 14).Value = 1.01701682785673

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032890107920676
$ws.Cells.Item(23, 4).Value = 1.040423362898899
$ws.Cells.Item(23, 5).Value = 1.032359263968911
$ws.Cells.Item(23, 6).Value = 1.048116232753378
$ws.Cells.Item(23, 9).Value = 1.037911006224966
$ws.Cells.Item(23, 10).Value = 1.039613491918788
$ws.Cells.Item(23, 11).Value = 1.044047513008483
$ws.Cells.Item(23, 12).Value = 1.036013943167879
$ws.Cells.Item(23, 13).Value = 1.051711743645609
$ws.Cells.Item(23, 14).Value = 1.017094139278891

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034310495504709
$ws.Cells.Item(24, 4).Value = 1.041547308627475
$ws.Cells.Item(24, 5).Value = 1.033563622047828
$ws.Cells.Item(24, 6).Value = 1.049580948812179
$ws.Cells.Item(24, 9).Value = 1.038280576707061
$ws.Cells.Item(24, 10).Value = 1.040511431472417
$ws.Cells.Item(24, 11).Value = 1.04489812529062
$ws.Cells.Item(24, 12).Value = 1.036942207696206
$ws.Cells.Item(24, 13).Value = 1.052904291595242
$ws.Cells.Item(24, 14).Value = 1.017398185802886

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035960740658548
$ws.Cells.Item(25, 4).Value = 1.042853433762743
$ws.Cells.Item(25, 5).Value = 1.034964617317231
$ws.Cells.Item(25, 6).Value = 1.051284932369591
$ws.Cells.Item(25, 9).Value = 1.03870495598843
$ws.Cells.Item(25, 10).Value = 1.041552578261882
$ws.Cells.Item(25, 11).Value = 1.045883959859936
$ws.Cells.Item(25, 12).Value = 1.038019774287592
$ws.Cells.Item(25, 13).Value = 1.054289579107341
$ws.Cells.Item(25, 14).Value = 1.017750336275501

